$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "42.613.71"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.255.41"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.52"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.25"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.05"
$ws.Range("E10").Value = "  +9.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0953"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.31"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.65"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.860"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "2.252.73"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "42.402.60"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.25"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.53"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.30"
$ws.Range("E23").Value = "  +26.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.60"
$ws.Range("E24").Value = "  +4.92%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.19"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0826"
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.73"
$ws.Range("E32").Value = "  -5.39%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.39"
$ws.Range("E34").Value = "  +10.33%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.57"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0317"
$ws.Range("E37").Value = "  +5.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.67"
$ws.Range("E38").Value = "  +7.31%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.36"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.07"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.86"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +5.99%  "
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.16"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("E51").Value = "  +0.71%  "
